$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Ending" dialogue rows (36-46) ---------------------------------

# Row 36
$ws.Range("A36").Value = 'Ending'
$ws.Range("B36").Value = '? ? ?'
$ws.Range("C36").Value = 'Well Detective, it seems you''ve gotten yourself in deep.'
$ws.Range("D36").Value = 100
$ws.Range("E36").Value = 'No'
$ws.Range("F36").Value = 'None'
$ws.Range("G36").Value = 'None'
$ws.Range("H36").Value = 'Game end.'

# Row 37
$ws.Range("A37").Value = 'Ending'
$ws.Range("B37").Value = 'Detective'
$ws.Range("C37").Value = 'Who''s there?'
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = '-'
$ws.Range("F37").Value = 'None'
$ws.Range("G37").Value = 'None'
$ws.Range("H37").Value = 'prev'

# Row 38
$ws.Range("A38").Value = 'Ending'
$ws.Range("B38").Value = 'Alien Don'
$ws.Range("C38").Value = 'The Alien Mafia sends its regards.'
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = '-'
$ws.Range("F38").Value = 'None'
$ws.Range("G38").Value = 'None'
$ws.Range("H38").Value = 'prev'

# Row 39
$ws.Range("A39").Value = 'Ending'
$ws.Range("B39").Value = 'Alien Don'
$ws.Range("C39").Value = 'And I welcome you to our humble tractor beam testing facility.'
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = '-'
$ws.Range("F39").Value = 'None'
$ws.Range("G39").Value = 'None'
$ws.Range("H39").Value = 'prev'

# Row 40
$ws.Range("A40").Value = 'Ending'
$ws.Range("B40").Value = 'Detective'
$ws.Range("C40").Value = 'The Alien Mafia? Tractor beam testing facility? Wait...'
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = '-'
$ws.Range("F40").Value = 'None'
$ws.Range("G40").Value = 'None'
$ws.Range("H40").Value = 'prev'

# Row 41 (only A/B/C populated)
$ws.Range("A41").Value = 'Ending'
$ws.Range("B41").Value = 'Detective'
$ws.Range("C41").Value = 'The chalk outline then, that was you?'

# Row 42 (only A/B/C populated)
$ws.Range("A42").Value = 'Ending'
$ws.Range("B42").Value = 'Alien Don'
$ws.Range("C42").Value = 'Yes, and the xenon balloon smuggler!'

# Row 43
$ws.Range("A43").Value = 'Ending'
$ws.Range("B43").Value = 'Detective'
$ws.Range("C43").Value = 'There are so many safety violations around this park!'
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = '-'
$ws.Range("F43").Value = 'None'
$ws.Range("G43").Value = 'None'
$ws.Range("H43").Value = 'prev'

# Row 44
$ws.Range("A44").Value = 'Ending'
$ws.Range("B44").Value = 'Alien Don'
$ws.Range("C44").Value = 'Huh? Wait-'
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = '-'
$ws.Range("F44").Value = 'None'
$ws.Range("G44").Value = 'None'
$ws.Range("H44").Value = 'prev'

# Row 45
$ws.Range("A45").Value = 'Ending'
$ws.Range("B45").Value = 'Detective'
$ws.Range("C45").Value = 'OSHA''s got your number, Glibohp.'
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = '-'
$ws.Range("F45").Value = 'None'
$ws.Range("G45").Value = 'None'
$ws.Range("H45").Value = 'prev'

# Row 46
$ws.Range("A46").Value = 'Ending'
$ws.Range("B46").Value = 'Detective'
$ws.Range("C46").Value = 'And after what I''ve seen today, the Alien Mafia''s going down.'
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = '-'
$ws.Range("F46").Value = 'None'
$ws.Range("G46").Value = 'None'
$ws.Range("H46").Value = 'prev'

# --- Row heights to match the rest of the sheet -------------------------
for ($r = 36; $r -le 46; $r++) {
    $ws.Rows.Item($r).RowHeight = 30.5
}

# --- Restore view state (frozen header, selection on new last row) ------
$ws.Range("C47").Select()
